{"js": "// Set bottom vertical alignment on every cell of the table's first two\n// (header) rows \u2014 adds <w:vAlign w:val=\"bottom\"/> to each of those\n// <w:tcPr> elements, matching the target diff.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst headerRowCount = Math.min(2, table.rows.items.length);\nconst headerRows = table.rows.items.slice(0, headerRowCount);\n\nfor (const row of headerRows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of headerRows) {\n  for (const cell of row.cells.items) {\n    cell.verticalAlignment = Word.VerticalAlignment.bottom;\n  }\n}\nawait context.sync();\n", "ps1": "# Set bottom vertical alignment on every cell of the table's first two\n# (header) rows \u2014 adds <w:vAlign w:val=\"bottom\"/> to each of those\n# <w:tcPr> elements, matching the target diff.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$wdCellAlignVerticalBottom = 3\n$headerRowCount = [Math]::Min(2, $table.Rows.Count)\n\nfor ($r = 1; $r -le $headerRowCount; $r++) {\n    $row = $table.Rows.Item($r)\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $row.Cells.Item($c).VerticalAlignment = $wdCellAlignVerticalBottom\n    }\n}\n"}
